# Northwind dashboard refresh: updated KPI/report numbers after improving
# the order-simulation determinism + WasShipped handling in the extractor.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: KPIs
# ---------------------------------------------------------------------
$wsKpi = $wb.Worksheets.Item(1)
$wsKpi.Cells.Item(2, 1).Value = 8340.884000000002
$wsKpi.Cells.Item(2, 4).Value = 39
$wsKpi.Cells.Item(2, 5).Value = 173.7684166666667
$wsKpi.Cells.Item(2, 6).Value = 0.9895833333333334

# ---------------------------------------------------------------------
# Sheet 2: Ventes Mensuelles (monthly sales) - columns C (TotalSales) and
# E (TotalQuantity) change for every existing row; D (NumOrders) unchanged.
# ---------------------------------------------------------------------
$wsVentes = $wb.Worksheets.Item(2)
$ventes = @(
    @(2, 707.735,           49),
    @(3, 508.5125000000001, 33),
    @(4, 1340.185,          92),
    @(5, 3608.8235,         197),
    @(6, 1239.578,          107),
    @(7, 936.05,            69)
)
foreach ($row in $ventes) {
    $r = $row[0]
    $wsVentes.Cells.Item($r, 3).Value = $row[1]
    $wsVentes.Cells.Item($r, 5).Value = $row[2]
}

# ---------------------------------------------------------------------
# Sheet 3: Par Catégorie - table grew from 2 data rows to 16 data rows.
# ---------------------------------------------------------------------
$wsCat = $wb.Worksheets.Item(3)
$categories = @(
    @("Jams, Preserves",          1728.4,             6,  39),
    @("Dried Fruit & Nuts",       1490.4875,          11, 77),
    @("Sauces",                   1086.825,           6,  37),
    @("Beverages",                755.6279999999999,  11, 59),
    @("Pasta",                    648.675,            6,  29),
    @("Dairy Products",           617.7,              4,  19),
    @("Condiments",                466.9,             5,  31),
    @("Baked Goods & Mixes",      444.821,            8,  45),
    @("Candy",                    296.4375,           4,  24),
    @("Grains",                   235.2,              5,  36),
    @("Canned Meat",              206.28,             4,  25),
    @("Oil",                      146.2475,           2,  8),
    @("Canned Fruit & Vegetables",81.59,              11, 63),
    @("Cereal",                   68,                 3,  19),
    @("Chips, Snacks",            45.45,              4,  27),
    @("Soups",                    22.2425,            2,  9)
)
$r = 2
foreach ($row in $categories) {
    $wsCat.Cells.Item($r, 1).Value = $row[0]
    $wsCat.Cells.Item($r, 2).Value = $row[1]
    $wsCat.Cells.Item($r, 3).Value = $row[2]
    $wsCat.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 4: Top Produits - table grew from 3 data rows to 20 data rows.
# ---------------------------------------------------------------------
$wsProd = $wb.Worksheets.Item(4)
$products = @(
    @("Northwind Traders Marmalade",        1227.15,  16, 3),
    @("Northwind Traders Curry Sauce",      780,       20, 3),
    @("Northwind Traders Walnuts",          645.1875,  32, 5),
    @("Northwind Traders Mozzarella",       617.7,     19, 4),
    @("Northwind Traders Boysenberry Spread",501.25,   23, 3),
    @("Northwind Traders Dried Pears",      403.5,     14, 3),
    @("Northwind Traders Coffee",           393.3,     9,  2),
    @("Northwind Traders Ravioli",          363.675,   21, 4),
    @("Northwind Traders Dried Apples",     302.1,     6,  1),
    @("Northwind Traders Chocolate",        296.4375,  24, 4),
    @("Northwind Traders Gnocchi",          285,       8,  2),
    @("Northwind Traders Cajun Seasoning",  284.9,     15, 3),
    @("Northwind Traders Long Grain Rice",  235.2,     36, 5),
    @("Northwind Traders Mustard",          182,       16, 3),
    @("Northwind Traders Scones",           171,       18, 3),
    @("Northwind Traders Tomato Sauce",     170,       10, 1),
    @("Northwind Traders Brownie Mix",      161.121,   14, 4),
    @("Northwind Traders Crab Meat",        160.08,    9,  2),
    @("Northwind Traders Olive Oil",        146.2475,  8,  2),
    @("Northwind Traders Hot Pepper Sauce", 136.825,   7,  2)
)
$r = 2
foreach ($row in $products) {
    $wsProd.Cells.Item($r, 1).Value = $row[0]
    $wsProd.Cells.Item($r, 2).Value = $row[1]
    $wsProd.Cells.Item($r, 3).Value = $row[2]
    $wsProd.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 5: Par Pays - only TotalSales (B2) changes.
# ---------------------------------------------------------------------
$wsPays = $wb.Worksheets.Item(5)
$wsPays.Cells.Item(2, 2).Value = 8340.884

# ---------------------------------------------------------------------
# Sheet 6: Employés - ranking / values changed for all 8 rows.
# ---------------------------------------------------------------------
$wsEmp = $wb.Worksheets.Item(6)
$employees = @(
    @("Nancy Freehafer",     1792.0455,          12, 6),
    @("Anne Hellung-Larsen", 1700.465,           10, 5),
    @("Jan Kotas",           1504.468,           6,  3),
    @("Mariya Sergienko",    1228.9555,          8,  4),
    @("Michael Neipper",     931.827,            4,  2),
    @("Andrew Cencini",      773.9504999999999,  4,  3),
    @("Laura Giussani",      213.86,             2,  1),
    @("Robert Zare",         195.3125,           2,  1)
)
$r = 2
foreach ($row in $employees) {
    $wsEmp.Cells.Item($r, 1).Value = $row[0]
    $wsEmp.Cells.Item($r, 2).Value = $row[1]
    $wsEmp.Cells.Item($r, 3).Value = $row[2]
    $wsEmp.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
